$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.446.49'
$ws.Range('E2').Value = '  +3.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.315.24'
$ws.Range('E3').Value = '  +2.19%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.02'
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.62'
$ws.Range('E6').Value = '  +5.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.537'
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('E9').Value = '  +7.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.01'
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0815'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.113'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.04'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.676.44'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.06'
$ws.Range('E15').Value = '  +2.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.293.60'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.815'
$ws.Range('E17').Value = '  +2.64%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.353.86'
$ws.Range('E18').Value = '  +3.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.51'
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0923'
$ws.Range('E20').Value = '  +1.96%  '
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.49'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.68'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('E24').Value = '  +5.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.64'
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.98'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.75'
$ws.Range('E28').Value = '  +4.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.87'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.18'
$ws.Range('E30').Value = '  +2.99%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.66'
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.99'
$ws.Range('E32').Value = '  +3.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.33'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.54'
$ws.Range('E35').Value = '  +7.33%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.89'
$ws.Range('E36').Value = '  +3.63%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0747'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').Value = '  -3.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.108'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.86'
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.116'
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.34'
$ws.Range('E42').Value = '  +9.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.94'
$ws.Range('E43').Value = '  +4.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.31'
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.983.25'
$ws.Range('E45').Value = '  +1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0290'
$ws.Range('E46').Value = '  +2.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.02'
$ws.Range('E47').Value = '  +4.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.84'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '56.16'
$ws.Range('E49').Value = '  +5.41%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.96'
$ws.Range('E50').Value = '  +15.39%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.56'
$ws.Range('E51').Value = '  +2.19%  '
